$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4136.0454
$ws.Range("I40").Value = 3583.1667
$ws.Range("J40").Value = 4343.375
$ws.Range("K40").Value = 3583.1667
$ws.Range("L40").Value = 4343.375
$ws.Range("M40").Value = -3408.1667
$ws.Range("N40").Value = -4693.375
$ws.Range("H41").Value = 567.8
$ws.Range("J41").Value = 801.8570999999999
$ws.Range("L41").Value = 801.8570999999999
$ws.Range("N41").Value = -1681.8571
$ws.Range("H100").Value = 7978.3105
$ws.Range("I100").Value = 2586.8
$ws.Range("K100").Value = 2586.8
$ws.Range("M100").Value = -2045.8
$ws.Range("H106").Value = 33249.082
$ws.Range("I106").Value = 36406.75
$ws.Range("K106").Value = 36406.75
$ws.Range("M106").Value = -35775.75
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H135").Value = 2954.7
$ws.Range("I135").Value = 2954.7
$ws.Range("K135").Value = 26592.3
$ws.Range("M135").Value = -24057.3
$ws.Range("H137").Value = 5962.75
$ws.Range("I137").Value = 4120.7
$ws.Range("J137").Value = 15173
$ws.Range("K137").Value = 12362.1
$ws.Range("L137").Value = 45519
$ws.Range("M137").Value = -9812.099999999999
$ws.Range("N137").Value = -50619
$ws.Range("H138").Value = 1114387.8
$ws.Range("I138").Value = 4499.5
$ws.Range("J138").Value = 1431498.8
$ws.Range("K138").Value = 13498.5
$ws.Range("L138").Value = 4294496.4
$ws.Range("M138").Value = -8358.5
$ws.Range("N138").Value = -4304776.4
$ws.Range("N133").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2712.5417
$ws.Range("I45").Value = 2205.1
$ws.Range("K45").Value = 2205.1
$ws.Range("M45").Value = -1828.1
$ws.Range("H61").Value = 4579.4707
$ws.Range("I61").Value = 4448.25
$ws.Range("K61").Value = 4448.25
$ws.Range("M61").Value = -4236.25
$ws.Range("H74").Value = 2194.9714
$ws.Range("I74").Value = 1767.6451
$ws.Range("K74").Value = 1767.6451
$ws.Range("M74").Value = -893.6451
$ws.Range("H77").Value = 2194.9714
$ws.Range("I77").Value = 1767.6451
$ws.Range("K77").Value = 8838.2255
$ws.Range("M77").Value = -4470.2255
$ws.Range("H94").Value = 49950
$ws.Range("J94").Value = 49950
$ws.Range("L94").Value = 49950
$ws.Range("N94").Value = -51752
$ws.Range("H122").Value = 1441.3214
$ws.Range("I122").Value = 1139.1818
$ws.Range("K122").Value = 3417.5454
$ws.Range("M122").Value = -967.5454
$ws.Range("H132").Value = 1527.0652
$ws.Range("I132").Value = 1412.4048
$ws.Range("J132").Value = 2731
$ws.Range("K132").Value = 4237.2144
$ws.Range("L132").Value = 8193
$ws.Range("M132").Value = -1707.2144
$ws.Range("N132").Value = -13253
$ws.Range("H136").Value = 4579.4707
$ws.Range("I136").Value = 4448.25
$ws.Range("K136").Value = 13344.75
$ws.Range("M136").Value = -10794.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 2900
$ws.Range("I10").Value = 2900
$ws.Range("K10").Value = 2900
$ws.Range("M10").Value = -2760
$ws.Range("H81").Value = 10390
$ws.Range("J81").Value = 10390
$ws.Range("L81").Value = 10390
$ws.Range("N81").Value = -12512
$ws.Range("H84").Value = 10390
$ws.Range("J84").Value = 10390
$ws.Range("L84").Value = 31170
$ws.Range("N84").Value = -41778
$ws.Range("H105").Value = 2466.889
$ws.Range("I105").Value = 914.9
$ws.Range("K105").Value = 914.9
$ws.Range("M105").Value = 832.1
$ws.Range("H134").Value = 14632.333
$ws.Range("I134").Value = 14632.333
$ws.Range("K134").Value = 43896.999
$ws.Range("M134").Value = -41361.999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 668.93335
$ws.Range("I16").Value = 559.56525
$ws.Range("K16").Value = 559.56525
$ws.Range("M16").Value = -272.56525
$ws.Range("H31").Value = 2496.4902
$ws.Range("J31").Value = 2853.258
$ws.Range("L31").Value = 2853.258
$ws.Range("N31").Value = -3443.258
$ws.Range("H34").Value = 2496.4902
$ws.Range("J34").Value = 2853.258
$ws.Range("L34").Value = 2853.258
$ws.Range("N34").Value = -3257.258
$ws.Range("H58").Value = 4255.1113
$ws.Range("I58").Value = 2169.9
$ws.Range("K58").Value = 2169.9
$ws.Range("M58").Value = -1966.9
$ws.Range("H113").Value = 668.93335
$ws.Range("I113").Value = 559.56525
$ws.Range("K113").Value = 559.56525
$ws.Range("M113").Value = 1610.43475
$ws.Range("H134").Value = 5837.2
$ws.Range("I134").Value = 5225.385
$ws.Range("K134").Value = 15676.155
$ws.Range("M134").Value = -13141.155
$ws.Range("H136").Value = 4255.1113
$ws.Range("I136").Value = 2169.9
$ws.Range("K136").Value = 6509.700000000001
$ws.Range("M136").Value = -3959.700000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 3854.8667
$ws.Range("I38").Value = 385.8
$ws.Range("J38").Value = 5589.4
$ws.Range("K38").Value = 1157.4
$ws.Range("L38").Value = 16768.2
$ws.Range("M38").Value = -810.4000000000001
$ws.Range("N38").Value = -17462.2
$ws.Range("H56").Value = 757375.0600000001
$ws.Range("I56").Value = 757375.0600000001
$ws.Range("K56").Value = 757375.0600000001
$ws.Range("M56").Value = -756845.0600000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 61933
$ws.Range("H97").Value = 848.9706
$ws.Range("I97").Value = 482.14285
$ws.Range("K97").Value = 482.14285
$ws.Range("M97").Value = 13.85714999999999
$ws.Range("H122").Value = 1237.88
$ws.Range("I122").Value = 1302.1666
$ws.Range("J122").Value = 1072.5714
$ws.Range("K122").Value = 3906.4998
$ws.Range("L122").Value = 3217.7142
$ws.Range("M122").Value = -1456.4998
$ws.Range("N122").Value = -8117.7142
$ws.Range("H126").Value = 7123
$ws.Range("I126").Value = 3910.7
$ws.Range("J126").Value = 12476.833
$ws.Range("K126").Value = 11732.1
$ws.Range("L126").Value = 37430.499
$ws.Range("M126").Value = -9262.099999999999
$ws.Range("N126").Value = -42370.499
$ws.Range("H132").Value = 2166.8696
$ws.Range("I132").Value = 2187.524
$ws.Range("J132").Value = 1950
$ws.Range("K132").Value = 6562.572
$ws.Range("L132").Value = 5850
$ws.Range("M132").Value = -4032.572
$ws.Range("N132").Value = -10910

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6169.1924
$ws.Range("I46").Value = 3686.875
$ws.Range("K46").Value = 3686.875
$ws.Range("M46").Value = -3498.875
$ws.Range("H132").Value = 4771.75
$ws.Range("I132").Value = 4738.176
$ws.Range("K132").Value = 14214.528
$ws.Range("M132").Value = -11684.528
$ws.Range("H136").Value = 3787.7354
$ws.Range("I136").Value = 3529.75
$ws.Range("K136").Value = 10589.25
$ws.Range("M136").Value = -8039.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1364.8846
$ws.Range("I122").Value = 1189.6666
$ws.Range("K122").Value = 3568.9998
$ws.Range("M122").Value = -1118.9998
$ws.Range("H136").Value = 1534.289
$ws.Range("I136").Value = 1325.1794
$ws.Range("J136").Value = 2893.5
$ws.Range("K136").Value = 3975.5382
$ws.Range("L136").Value = 8680.5
$ws.Range("M136").Value = -1425.5382
$ws.Range("N136").Value = -13780.5
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
